$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()

# Clear values in A2:A285 but keep formatting
$ws.Range("A2:A285").ClearContents()

# Update frozen pane top-left cell
$ws.Range("O2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Set selection to a range
$ws.Range("A1:AG1").Select()
